$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 234; this shifts the existing rows 234:249
# down to 235:250 (and the sheet dimension grows to A1:R250 automatically).
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new weekly price record.
$ws.Cells.Item(234, 1).Value = 10
$ws.Cells.Item(234, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(234, 3).Value = "La Araucanía"
$ws.Cells.Item(234, 4).Value = 44753
$ws.Cells.Item(234, 5).Value = 9
$ws.Cells.Item(234, 6).Value = 100112043
$ws.Cells.Item(234, 7).Value = "Pepino dulce"
$ws.Cells.Item(234, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 180
$ws.Cells.Item(234, 11).Value = 18000
$ws.Cells.Item(234, 12).Value = 19000
$ws.Cells.Item(234, 13).Value = 18556
$ws.Cells.Item(234, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(234, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(234, 16).Value = 1031
$ws.Cells.Item(234, 17).Value = 18
$ws.Cells.Item(234, 18).Value = "Hortaliza"
